$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet from "UniformA-HW15.xpc" to "UniformA"
$ws.Name = "UniformA"

# Append a new data row (row 16) that follows the same pattern as row 15:
#  - Column A: next HKL group index (14), with the bordered/bold "header" style
#  - Column B: same label as row 15 (HexGrid-60degTilt5degRes)
#  - Columns C:P: intensity values of 1

$ws.Range("A15").Copy($ws.Range("A16"))
$ws.Range("A16").Value2 = 14

$ws.Range("B15").Copy($ws.Range("B16"))

for ($col = 3; $col -le 16; $col++) {
    $ws.Cells.Item(16, $col).Value2 = 1
}
